$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 3 for columns A, B, E, F, G, H, Q, R
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr2 = "${col}2"
    $addr3 = "${col}3"
    $cell2 = $ws.Range($addr2)
    $cell3 = $ws.Range($addr3)
    $tmp = $cell2.Value2
    $cell2.Value2 = $cell3.Value2
    $cell3.Value2 = $tmp
}
